$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '27.268.74'
$c.Style = "Normal"
$ws.Range('E2').Value = '  -0.69%  '
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '1.702.42'
$c.Style = "Normal"
$ws.Range('E3').Value = '  -1.23%  '
$c = $ws.Range('D4')
$c.NumberFormat = "@"
$c.Value = '1.003'
$c.Style = "Normal"
$ws.Range('E4').Value = '  -0.17%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '223.37'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -1.01%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '0.5304'
$c.Style = "Normal"
$ws.Range('E6').Value = '  -1.14%  '
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.Style = "Normal"
$ws.Range('E7').Value = '  -0.23%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.2657'
$c.Style = "Normal"
$ws.Range('E8').Value = '  -0.81%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.06575'
$c.Style = "Normal"
$ws.Range('E9').Value = '  -0.21%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '20.70'
$c.Style = "Normal"
$ws.Range('E10').Value = '  -4.24%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.07619'
$c.Style = "Normal"
$ws.Range('E11').Value = '  -1.85%  '
$ws.Range('E12').Value = '  -3.09%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '1.717.40'
$c.Style = "Normal"
$ws.Range('E13').Value = '  -0.11%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '1.937.98'
$c.Style = "Normal"
$ws.Range('E14').Value = '  -1.09%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '0.5772'
$c.Style = "Normal"
$ws.Range('E15').Value = '  -1.62%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '0.0₅8142'
$c.Style = "Normal"
$ws.Range('E16').Value = '  -1.37%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '67.43'
$c.Style = "Normal"
$ws.Range('E17').Value = '  -0.71%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '27.268.46'
$c.Style = "Normal"
$ws.Range('E18').Value = '  -0.78%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '215.16'
$c.Style = "Normal"
$ws.Range('E19').Value = '  -3.45%  '
$ws.Range('E20').Value = '  -0.19%  '
$ws.Range('E21').Value = '  -2.71%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '10.34'
$c.Style = "Normal"
$ws.Range('E22').Value = '  -3.14%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '5.958'
$c.Style = "Normal"
$ws.Range('E23').Value = '  -2.29%  '
$ws.Range('E24').Value = '  -0.34%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '143.76'
$c.Style = "Normal"
$ws.Range('E25').Value = '  -2.84%  '
$ws.Range('E26').Value = '  +0.78%  '
$ws.Range('E27').Value = '  -2.66%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '7.195'
$c.Style = "Normal"
$ws.Range('E28').Value = '  -2.84%  '
$ws.Range('E29').Value = '  -3.51%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '0.05355'
$c.Style = "Normal"
$ws.Range('E30').Value = '  -3.59%  '
$ws.Range('E31').Value = '  -1.50%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '3.456'
$c.Style = "Normal"
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '3.394'
$c.Style = "Normal"
$ws.Range('E33').Value = '  -1.88%  '
$ws.Range('E34').Value = '  -1.17%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '2.865'
$c.Style = "Normal"
$ws.Range('E35').Value = '  +1.73%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '2.411'
$c.Style = "Normal"
$ws.Range('E36').Value = '  -1.65%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '0.9446'
$c.Style = "Normal"
$ws.Range('E37').Value = '  -1.41%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '0.5797'
$c.Style = "Normal"
$ws.Range('E38').Value = '  -2.03%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '0.01625'
$c.Style = "Normal"
$ws.Range('E39').Value = '  -1.04%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '5.776'
$c.Style = "Normal"
$ws.Range('E40').Value = '  -1.54%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.Style = "Normal"
$ws.Range('E41').Value = '  -0.29%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '1.039.15'
$c.Style = "Normal"
$ws.Range('E42').Value = '  -1.51%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '0.8388'
$c.Style = "Normal"
$ws.Range('E43').Value = '  -1.98%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '100.83'
$c.Style = "Normal"
$ws.Range('E44').Value = '  -0.68%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '1.845.32'
$c.Style = "Normal"
$ws.Range('E45').Value = '  -1.09%  '
$ws.Range('E46').Value = '  -0.30%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '57.69'
$c.Style = "Normal"
$ws.Range('E47').Value = '  -2.06%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '0.4518'
$c.Style = "Normal"
$ws.Range('E48').Value = '  +1.61%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '1.004'
$c.Style = "Normal"
$ws.Range('E49').Value = '  +0.36%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '8.028'
$c.Style = "Normal"
$ws.Range('E50').Value = '  -1.98%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '0.05227'
$c.Style = "Normal"
$ws.Range('E51').Value = '  -0.97%  '
